# Homework from the second seminar: append the "Урок 2" section
# (heading + assignment text + surrounding blank paragraphs) to the
# end of the document, right before the section break.

$d = $word.ActiveDocument

# Collapse a range at the very end of the document content (right
# before the final section mark) so the new paragraphs land after the
# existing "... кто вытянул последнюю спичку." paragraph.
$r = $d.Content
$r.Collapse(0)

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

$newParagraphsXml = @"
<w:p $wordNs>
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/>
      <w:color w:val="2C2D30"/>
      <w:kern w:val="0"/>
      <w:sz w:val="23"/>
      <w:szCs w:val="23"/>
      <w:lang w:val="en-US" w:eastAsia="ru-RU"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $wordNs>
  <w:pPr>
    <w:pStyle w:val="3"/>
    <w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0" w:line="360" w:lineRule="atLeast"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="3F5368"/>
      <w:sz w:val="26"/>
      <w:szCs w:val="26"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="3F5368"/>
      <w:sz w:val="26"/>
      <w:szCs w:val="26"/>
    </w:rPr>
    <w:t>Урок 2. Библиотеки. Статические и динамические</w:t>
  </w:r>
</w:p>
<w:p $wordNs>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
      <w:color w:val="2C2D30"/>
      <w:sz w:val="23"/>
      <w:szCs w:val="23"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/>
      <w:color w:val="2C2D30"/>
      <w:sz w:val="23"/>
      <w:szCs w:val="23"/>
    </w:rPr>
    <w:t>Задание: реализуйте классы для игры в спички из предыдущего занятия в виде библиотеки, после чего используйте ее для реализации исполняемого файла с помощью собранной библиотеки.</w:t>
  </w:r>
</w:p>
<w:p $wordNs>
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/>
      <w:color w:val="2C2D30"/>
      <w:kern w:val="0"/>
      <w:sz w:val="23"/>
      <w:szCs w:val="23"/>
      <w:lang w:val="en-US" w:eastAsia="ru-RU"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $wordNs>
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/>
      <w:color w:val="2C2D30"/>
      <w:kern w:val="0"/>
      <w:sz w:val="23"/>
      <w:szCs w:val="23"/>
      <w:lang w:val="en-US" w:eastAsia="ru-RU"/>
      <w14:ligatures w14:val="none"/>
    </w:rPr>
  </w:pPr>
</w:p>
"@

$r.InsertXML($newParagraphsXml)
